# Update version string across the workbook:
#   old: mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)
#   new: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: standalone version string
$aboutSheet.Range("A2").Value2 = "Version: " + $newVersion

# A6: citation text embedding the version string
$oldCitation = [string]$aboutSheet.Range("A6").Value2
$newCitation = $oldCitation -replace [regex]::Escape($oldVersion), $newVersion
$aboutSheet.Range("A6").Value2 = $newCitation

# S2:S8 on the data sheet hold the build_version value for each row
for ($r = 2; $r -le 8; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)  # column S = 19
    if ([string]$cell.Value2 -eq $oldVersion) {
        $cell.Value2 = $newVersion
    }
}
